# 2021 Team Data.xlsx -- log 2021 divisional round, simulated season from
# conference round.
#
# The workbook tracks running per-play logs (space separated numbers held
# in single cells) plus season-to-date summary totals on each tab. This
# script appends one more game's worth of data to each running log and
# bumps every summary total that depends on it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS tab - per-play rush/pass yards logs for OFF (col B) and DEF (col C)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 2 -1 8 3 -1 2 4 3 4 6 3 8 9 3 9 45 5 5 1 -1 9 10 0 3 0"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 1 1 2 5 -3 0 1 2 6 10 7 16 1 13 -1 4 2 -2"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 11 13 3 41 3 13 8 40 5 33 20 6 3 16 5"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 21 8 4 57 7 19 7 7 6 15 4 22 1 5 10 11 32 7 12 11 4 7 21 8 13 7 3 19"

# ---------------------------------------------------------------------
# OFF tab - season totals, Home row (r2) and Road row (r3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 271
$ws.Range("F2").Value = 95
$ws.Range("G2").Value = 66
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 40
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = 202
$ws.Range("O2").Value = 25
$ws.Range("Q2").Value = 622

$ws.Range("C3").Value = 154
$ws.Range("E3").Value = 44
$ws.Range("F3").Value = 114
$ws.Range("G3").Value = 26
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 67
$ws.Range("N3").Value = 24

# ---------------------------------------------------------------------
# DEF tab - season totals, Home row (r2) and Road row (r3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 172
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 53
$ws.Range("G2").Value = 63
$ws.Range("J2").Value = 28
$ws.Range("L2").Value = 359
$ws.Range("M2").Value = 224
$ws.Range("Q2").Value = 607

$ws.Range("C3").Value = 222
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 126
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 68
$ws.Range("J3").Value = 61
$ws.Range("N3").Value = 31

# ---------------------------------------------------------------------
# ST tab - kicking/punting totals + per-game logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 94
$ws.Range("D2").Value = 61
$ws.Range("F2").Value = 271
$ws.Range("G2").Value = 258
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 118

$ws.Range("B3").Value = 53

$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 61 57 64"
$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 31 32 19"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 26"
$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 50 27 41 48"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 0 0 0 0"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 5 0 0 0"

# ---------------------------------------------------------------------
# TURNS tab
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 11
$ws.Range("D3").Value = 8

# ---------------------------------------------------------------------
# PEN tab
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value = 17
